# Apply row permutation described in the commit diff.
# Rows 2,3,5,6,8,9,10 get the data that used to live in a different row
# (rows 4 and 7 are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <= original row 6 data
$ws.Range("A2").Value = 111609175
$ws.Range("B2").Value = 77268
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 228912
$ws.Range("F2").Value = "Mörk kolflarnlav"
$ws.Range("G2").Value = "Carbonicola myrmecina"
$ws.Range("H2").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("I2").Value = ""
$ws.Range("Q2").Value = 514769.8196280882
$ws.Range("R2").Value = 6925156.6384242
$ws.Range("AO2").Value = "brandstubbe"

# Row 3 <= original row 10 data
$ws.Range("A3").Value = 111609169
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("I3").Value = "'4"
$ws.Range("Q3").Value = 515078.8479096842
$ws.Range("R3").Value = 6925177.45879681
$ws.Range("AO3").Value = ""

# Row 5 <= original row 8 data
$ws.Range("A5").Value = 111609168
$ws.Range("B5").Value = 77597
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 864
$ws.Range("F5").Value = "Knottrig blåslav"
$ws.Range("G5").Value = "Hypogymnia bitteri"
$ws.Range("H5").Value = "(Lynge) Ahti"
$ws.Range("I5").Value = ""
$ws.Range("Q5").Value = 515085.0087401169
$ws.Range("R5").Value = 6925147.4056778
$ws.Range("AO5").Value = "tall"

# Row 6 <= original row 9 data
$ws.Range("A6").Value = 111609173
$ws.Range("B6").Value = 96348
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("I6").Value = "'7"
$ws.Range("Q6").Value = 514934.1293421969
$ws.Range("R6").Value = 6925308.234934391
$ws.Range("AO6").Value = ""

# Row 8 <= original row 3 data
$ws.Range("A8").Value = 111609172
$ws.Range("B8").Value = 77268
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 228912
$ws.Range("F8").Value = "Mörk kolflarnlav"
$ws.Range("G8").Value = "Carbonicola myrmecina"
$ws.Range("H8").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("I8").Value = ""
$ws.Range("Q8").Value = 514955.9350709137
$ws.Range("R8").Value = 6925302.779521272
$ws.Range("AO8").Value = "brandstubbe"

# Row 9 <= original row 2 data
$ws.Range("A9").Value = 111609167
$ws.Range("B9").Value = 77186
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 353
$ws.Range("F9").Value = "Dvärgbägarlav"
$ws.Range("G9").Value = "Cladonia parasitica"
$ws.Range("H9").Value = "(Hoffm.) Hoffm."
$ws.Range("I9").Value = ""
$ws.Range("Q9").Value = 515051.1877758073
$ws.Range("R9").Value = 6925144.938876954
$ws.Range("AO9").Value = "silverlåga av tall"

# Row 10 <= original row 5 data
$ws.Range("A10").Value = 111609170
$ws.Range("B10").Value = 96348
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = "Knärot"
$ws.Range("G10").Value = "Goodyera repens"
$ws.Range("H10").Value = "(L.) R. Br."
$ws.Range("I10").Value = "'3"
$ws.Range("Q10").Value = 515035.9338400747
$ws.Range("R10").Value = 6925238.814452391
$ws.Range("AO10").Value = ""
